$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 1) "Test Writer:" paragraph -> split into multiple runs with proofErr marks
#    around the two names that Word's spell-checker would flag.
# ---------------------------------------------------------------------------
$cellWriter = $tbl.Cell(1, 1)
$pWriter = $cellWriter.Range.Paragraphs.Item(1)
$rWriter = $pWriter.Range

$writerXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w:rsidR="00523A40" w:rsidRDefault="00523A40" w:rsidP="005F2733">' +
  '<w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
  '<w:r><w:t>Test Writer</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Edgard</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Musafiri</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> and Sarmad Butti</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$rWriter.InsertXML($writerXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Description cell text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
  "Testing the matrix of LEDs and making sure they are all functional.",
  $false, $false, $false, $false, $false, $true, 1, $false,
  "Ensuring all LEDs in the LED matrix are functioning. Observing how the LED matrix responds.",
  2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Setup cell text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
  "Attaching a battery to the device and start testing ",
  $false, $false, $false, $false, $false, $true, 1, $false,
  "Connecting the device to a battery and begin testing the LED response.",
  2) | Out-Null

# ---------------------------------------------------------------------------
# 4) LED 5-9 "turn ON" cells: merge the two runs into a single run, and
#    relocate the _GoBack bookmark that used to trail the LED 9 cell.
# ---------------------------------------------------------------------------
foreach ($rowIdx in 13..17) {
  $ledNum = $rowIdx - 8
  $cell = $tbl.Rows.Item($rowIdx).Cells.Item(3)
  $p = $cell.Range.Paragraphs.Item(1)
  $r = $p.Range
  $r.Delete()
  $r.InsertAfter("LED $ledNum turn ON")
}

# ---------------------------------------------------------------------------
# 5) Hardware Ver cell: replace the literal "1" with an empty paragraph that
#    now hosts the relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$cellHwVer = $tbl.Rows.Item(6).Cells.Item(2)
$pHwVer = $cellHwVer.Range.Paragraphs.Item(1)
$rHwVer = $pHwVer.Range
$rHwVer.Delete()
$d.Bookmarks.Add("_GoBack", $rHwVer) | Out-Null

Write-Output "done"
